$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $c = $ws.Range($cell)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = $origStyle
}

Set-TextValue 'D2' '45.124.94'
Set-TextValue 'E2' '  -3.30%  '
Set-TextValue 'D3' '2.387.58'
Set-TextValue 'E3' '  +5.77%  '
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '292.50'
Set-TextValue 'D6' '93.37'
Set-TextValue 'E6' '  -6.79%  '
Set-TextValue 'E7' '  -0.74%  '
Set-TextValue 'E8' '  +0.02%  '
Set-TextValue 'E9' '  -2.41%  '
Set-TextValue 'D10' '34.11'
Set-TextValue 'E10' '  -3.91%  '
Set-TextValue 'D11' '0.0774'
Set-TextValue 'E11' '  -0.68%  '
Set-TextValue 'E12' '  -2.72%  '
Set-TextValue 'E13' '  +0.94%  '
Set-TextValue 'D14' '2.753.18'
Set-TextValue 'E14' '  +5.75%  '
Set-TextValue 'D15' '2.388.24'
Set-TextValue 'E15' '  +5.33%  '
Set-TextValue 'D16' '14.07'
Set-TextValue 'E16' '  +3.99%  '
Set-TextValue 'D17' '0.825'
Set-TextValue 'E17' '  +3.21%  '
Set-TextValue 'D18' '45.133.40'
Set-TextValue 'E18' '  -3.24%  '
Set-TextValue 'D19' '12.37'
Set-TextValue 'E19' '  -4.26%  '
Set-TextValue 'D20' '0.0₃0933'
Set-TextValue 'E20' '  +0.75%  '
Set-TextValue 'D21' '6.07'
Set-TextValue 'E21' '  +3.19%  '
Set-TextValue 'D22' '66.39'
Set-TextValue 'E22' '  +1.76%  '
Set-TextValue 'D23' '237.39'
Set-TextValue 'E23' '  -4.63%  '
Set-TextValue 'E24' '  -2.85%  '
Set-TextValue 'E25' '  +0.03%  '
Set-TextValue 'E26' '  +0.85%  '
Set-TextValue 'D27' '2.21'
Set-TextValue 'E27' '  -1.11%  '
Set-TextValue 'D28' '37.28'
Set-TextValue 'E28' '  -12.65%  '
Set-TextValue 'D29' '9.51'
Set-TextValue 'E29' '  -1.76%  '
Set-TextValue 'D30' '3.81'
Set-TextValue 'E30' '  +19.28%  '
Set-TextValue 'D31' '20.84'
Set-TextValue 'E31' '  +5.09%  '
Set-TextValue 'D32' '2.71'
Set-TextValue 'E32' '  -2.73%  '
Set-TextValue 'D33' '146.60'
Set-TextValue 'D34' '5.37'
Set-TextValue 'E34' '  -0.85%  '
Set-TextValue 'D35' '0.0755'
Set-TextValue 'E35' '  -1.70%  '
Set-TextValue 'E36' '  +14.43%  '
Set-TextValue 'D37' '0.111'
Set-TextValue 'E37' '  -1.51%  '
Set-TextValue 'E38' '  -1.29%  '
Set-TextValue 'D39' '14.42'
Set-TextValue 'E39' '  -10.86%  '
Set-TextValue 'E40' '  -4.39%  '
Set-TextValue 'E41' '  -1.70%  '
Set-TextValue 'D42' '1.970.16'
Set-TextValue 'E42' '  +8.62%  '
Set-TextValue 'D43' '3.15'
Set-TextValue 'E43' '  -1.29%  '
Set-TextValue 'E44' '  -0.08%  '
Set-TextValue 'D45' '87.96'
Set-TextValue 'E45' '  -2.42%  '
Set-TextValue 'D46' '1.69'
Set-TextValue 'E46' '  -14.02%  '
Set-TextValue 'D47' '15.45'
Set-TextValue 'E47' '  +20.72%  '
Set-TextValue 'D48' '8.38'
Set-TextValue 'E48' '  +8.01%  '
Set-TextValue 'D49' '98.92'
Set-TextValue 'E49' '  +5.70%  '
Set-TextValue 'D50' '2.624.39'
Set-TextValue 'E50' '  +5.81%  '
Set-TextValue 'D51' '0.181'
Set-TextValue 'E51' '  -3.55%  '
